$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental = true  (row 7, column B)
$ws.Range("B7").Formula = '="true"'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date updated (row 8, column B)
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"

# Compositional = false (row 18, column B)
$ws.Range("B18").Formula = '="false"'
$ws.Range("B18").Copy()
$ws.Range("B18").PasteSpecial(-4163)
